# Auto-generated edit script: update crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.489.01"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "2.989.40"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'381.74"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "'103.50"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("D7").Value = "'0.547"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.594"
$ws.Range("D10").Value = "'36.75"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").Value = "'0.0862"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").Value = "3.460.95"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").Value = "'7.82"
$ws.Range("E14").Value = "  +3.90%  "
$ws.Range("D15").Value = "'18.45"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").Value = "3.007.39"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("D17").Value = "'11.33"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "'1.00"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").Value = "51.516.30"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "'3.15"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'12.63"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").Value = "'70.36"
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("D24").Value = "'268.05"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("E25").Value = "  +3.18%  "
$ws.Range("D26").Value = "'7.85"
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("D27").Value = "'7.52"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'26.14"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").Value = "'0.166"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").Value = "'0.110"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").Value = "'10.36"
$ws.Range("E32").Value = "  +3.90%  "
$ws.Range("D33").Value = "'34.79"
$ws.Range("E33").Value = "  +4.47%  "
$ws.Range("D34").Value = "'51.50"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").Value = "'0.0443"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'3.26"
$ws.Range("E38").Value = "  +2.86%  "
$ws.Range("D39").Value = "'16.82"
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.85"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.56"
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").Value = "'124.86"
$ws.Range("E43").Value = "  +4.07%  "
$ws.Range("D44").Value = "'3.66"
$ws.Range("E44").Value = "  +9.64%  "
$ws.Range("D45").Value = "'21.71"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.273"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'2.03"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("D49").Value = "2.038.10"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("E51").Value = "  +16.16%  "
